# Weekly update: insert a new daily price record for
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Pepino ensalada".
#
# The new record is inserted as row 133 (the most recent week's data point),
# which pushes the previously existing rows 133-153 down to rows 134-154.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 133, shifting rows 133:153 down to 134:154.
$ws.Rows("133:133").Insert()

# Populate the newly inserted row 133 with the new weekly data point.
$ws.Range("A133").Value = 7
$ws.Range("B133").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C133").Value = "Ñuble"
$ws.Range("D133").Value = 44476
$ws.Range("E133").Value = 16
$ws.Range("F133").Value = 100112043
$ws.Range("G133").Value = "Pepino ensalada"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 120
$ws.Range("K133").Value = 17000
$ws.Range("L133").Value = 18000
$ws.Range("M133").Value = 17500
$ws.Range("N133").Value = "$/caja 60 unidades"
$ws.Range("O133").Value = "Región de Arica y Parinacota"
$ws.Range("P133").Value = 292
$ws.Range("Q133").Value = 60
$ws.Range("R133").Value = "Hortaliza"
